# "fixed export and fixing maps"
#
# The visible data sheet (currently named "1") is renamed to "ხობი" and its
# small summary table is trimmed down: the extra "(census results)" note
# row is removed, and the historical 1989 / 2002 area columns are dropped,
# leaving only the 2014 figure next to the area ("ფართობი") row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the data sheet from "1" to "ხობი".
$ws.Name = "ხობი"

# Remove the "(census results)" note row (old row 2) - everything below
# shifts up one row.
$ws.Rows.Item(2).Delete()

# Drop the 1989 and 2002 columns, keeping only the 2014 figure (which
# shifts from column D into column B).
$ws.Columns.Item(2).Delete()
$ws.Columns.Item(2).Delete()

# Match the saved selection state (active cell A2 on the trimmed sheet).
[void]$ws.Range("A2").Select()
